# Refresh the LR-pair edge table (rows 2-10) with the recomputed TPM-based
# NATMI statistics -- now covering the full ECs/FAPs/MuSCs sending x target
# cluster cross-product (9 rows) instead of the previous partial 6-row table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnfsf13"
$ws.Cells.Item(2, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.6848073333333332
$ws.Cells.Item(2, 8).Value = 2.054422
$ws.Cells.Item(2, 9).Value = 0.2268310526442471
$ws.Cells.Item(2, 10).Value = 0.2268310526442472
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.07580833333333332
$ws.Cells.Item(2, 14).Value = 0.227425
$ws.Cells.Item(2, 15).Value = 0.028190957994264
$ws.Cells.Item(2, 16).Value = 0.02819095799426401
$ws.Cells.Item(2, 17).Value = 0.05191410259444443
$ws.Cells.Item(2, 18).Value = 0.4672269233499999
$ws.Cells.Item(2, 19).Value = 0.006394584676888657
$ws.Cells.Item(2, 20).Value = 0.00639458467688866

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnfsf13"
$ws.Cells.Item(3, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.6848073333333332
$ws.Cells.Item(3, 8).Value = 2.054422
$ws.Cells.Item(3, 9).Value = 0.2268310526442471
$ws.Cells.Item(3, 10).Value = 0.2268310526442472
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 2.095195666666667
$ws.Cells.Item(3, 14).Value = 6.285587
$ws.Cells.Item(3, 15).Value = 0.7791435378093522
$ws.Cells.Item(3, 16).Value = 0.7791435378093522
$ws.Cells.Item(3, 17).Value = 1.434805357301555
$ws.Cells.Item(3, 18).Value = 12.913248215714
$ws.Cells.Item(3, 19).Value = 0.1767339488422581
$ws.Cells.Item(3, 20).Value = 0.1767339488422582

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnfsf13"
$ws.Cells.Item(4, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.6848073333333332
$ws.Cells.Item(4, 8).Value = 2.054422
$ws.Cells.Item(4, 9).Value = 0.2268310526442471
$ws.Cells.Item(4, 10).Value = 0.2268310526442472
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.5180969999999999
$ws.Cells.Item(4, 14).Value = 1.554291
$ws.Cells.Item(4, 15).Value = 0.1926655041963838
$ws.Cells.Item(4, 16).Value = 0.1926655041963838
$ws.Cells.Item(4, 17).Value = 0.3547966249779999
$ws.Cells.Item(4, 18).Value = 3.193169624801999
$ws.Cells.Item(4, 19).Value = 0.04370251912510035
$ws.Cells.Item(4, 20).Value = 0.04370251912510036

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnfsf13"
$ws.Cells.Item(5, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.3279213333333333
$ws.Cells.Item(5, 8).Value = 0.983764
$ws.Cells.Item(5, 9).Value = 0.1086184939966157
$ws.Cells.Item(5, 10).Value = 0.1086184939966157
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.07580833333333332
$ws.Cells.Item(5, 14).Value = 0.227425
$ws.Cells.Item(5, 15).Value = 0.028190957994264
$ws.Cells.Item(5, 16).Value = 0.02819095799426401
$ws.Cells.Item(5, 17).Value = 0.02485916974444444
$ws.Cells.Item(5, 18).Value = 0.2237325277
$ws.Cells.Item(5, 19).Value = 0.003062059401658809
$ws.Cells.Item(5, 20).Value = 0.00306205940165881

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnfsf13"
$ws.Cells.Item(6, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.3279213333333333
$ws.Cells.Item(6, 8).Value = 0.983764
$ws.Cells.Item(6, 9).Value = 0.1086184939966157
$ws.Cells.Item(6, 10).Value = 0.1086184939966157
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.095195666666667
$ws.Cells.Item(6, 14).Value = 6.285587
$ws.Cells.Item(6, 15).Value = 0.7791435378093522
$ws.Cells.Item(6, 16).Value = 0.7791435378093522
$ws.Cells.Item(6, 17).Value = 0.6870593566075556
$ws.Cells.Item(6, 18).Value = 6.183534209467999
$ws.Cells.Item(6, 19).Value = 0.08462939768404702
$ws.Cells.Item(6, 20).Value = 0.08462939768404702

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnfsf13"
$ws.Cells.Item(7, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.3279213333333333
$ws.Cells.Item(7, 8).Value = 0.983764
$ws.Cells.Item(7, 9).Value = 0.1086184939966157
$ws.Cells.Item(7, 10).Value = 0.1086184939966157
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.5180969999999999
$ws.Cells.Item(7, 14).Value = 1.554291
$ws.Cells.Item(7, 15).Value = 0.1926655041963838
$ws.Cells.Item(7, 16).Value = 0.1926655041963838
$ws.Cells.Item(7, 17).Value = 0.169895059036
$ws.Cells.Item(7, 18).Value = 1.529055531324
$ws.Cells.Item(7, 19).Value = 0.02092703691090985
$ws.Cells.Item(7, 20).Value = 0.02092703691090985

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Tnfsf13"
$ws.Cells.Item(8, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.006290666666667
$ws.Cells.Item(8, 8).Value = 6.018872
$ws.Cells.Item(8, 9).Value = 0.6645504533591371
$ws.Cells.Item(8, 10).Value = 0.6645504533591372
$ws.Cells.Item(8, 11).Value = 1
$ws.Cells.Item(8, 12).Value = 0.3333333333333333
$ws.Cells.Item(8, 13).Value = 0.07580833333333332
$ws.Cells.Item(8, 14).Value = 0.227425
$ws.Cells.Item(8, 15).Value = 0.028190957994264
$ws.Cells.Item(8, 16).Value = 0.02819095799426401
$ws.Cells.Item(8, 17).Value = 0.1520935516222222
$ws.Cells.Item(8, 18).Value = 1.3688419646
$ws.Cells.Item(8, 19).Value = 0.01873431391571653
$ws.Cells.Item(8, 20).Value = 0.01873431391571654

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Tnfsf13"
$ws.Cells.Item(9, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.006290666666667
$ws.Cells.Item(9, 8).Value = 6.018872
$ws.Cells.Item(9, 9).Value = 0.6645504533591371
$ws.Cells.Item(9, 10).Value = 0.6645504533591372
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.095195666666667
$ws.Cells.Item(9, 14).Value = 6.285587
$ws.Cells.Item(9, 15).Value = 0.7791435378093522
$ws.Cells.Item(9, 16).Value = 0.7791435378093522
$ws.Cells.Item(9, 17).Value = 4.203571510873778
$ws.Cells.Item(9, 18).Value = 37.832143597864
$ws.Cells.Item(9, 19).Value = 0.517780191283047
$ws.Cells.Item(9, 20).Value = 0.5177801912830471

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Tnfsf13"
$ws.Cells.Item(10, 3).Value = "Tnfrsf11b"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.006290666666667
$ws.Cells.Item(10, 8).Value = 6.018872
$ws.Cells.Item(10, 9).Value = 0.6645504533591371
$ws.Cells.Item(10, 10).Value = 0.6645504533591372
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 0.5180969999999999
$ws.Cells.Item(10, 14).Value = 1.554291
$ws.Cells.Item(10, 15).Value = 0.1926655041963838
$ws.Cells.Item(10, 16).Value = 0.1926655041963838
$ws.Cells.Item(10, 17).Value = 1.039453175528
$ws.Cells.Item(10, 18).Value = 9.355078579752
$ws.Cells.Item(10, 19).Value = 0.1280359481603736
$ws.Cells.Item(10, 20).Value = 0.1280359481603736

